# Commit: "update fig1map and fig1"
# The underlying data table (sheet1) dropped four rows that referenced
# studies/locations which are no longer part of the figure map:
#   - DNK1  / ageband   (row 6)
#   - BRA4  / region    (row 19)
#   - SF_CA1/ region    (row 22)
#   - DNK1_nch / ageband (row 27)
#
# Deleting whole rows (rather than just clearing cell contents) shifts
# everything below up, which is exactly what the target workbook shows
# (dimension shrinks from A1:E28 to A1:E24, and all of the subsequent
# rows move up by the number of preceding deleted rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(6).Delete()

# Restore the active sheet/selection shown in the edited workbook.
$ws.Activate()
$ws.Range("F21").Select()
